# Weekly fruit/vegetable price update: insert two new daily-price records
# for "Alcachofa" (Vega Modelo de Temuco) right before the existing row 200,
# shifting the rest of the data table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 200 (old rows 200-293 shift down to 202-295).
$ws.Range("200:201").Insert()

# --- New row 200 ---
$ws.Range("A200").Value = 10
$ws.Range("B200").Value = "Vega Modelo de Temuco"
$ws.Range("C200").Value = "La Araucanía"
$ws.Range("D200").Value = 45089
$ws.Range("E200").Value = 9
$ws.Range("F200").Value = 100112013
$ws.Range("G200").Value = "Alcachofa"
$ws.Range("H200").Value = "Española"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 35
$ws.Range("K200").Value = 20000
$ws.Range("L200").Value = 20000
$ws.Range("M200").Value = 20000
$ws.Range("N200").Value = "$/caja 35 unidades"
$ws.Range("O200").Value = "Provincia de Limarí"
$ws.Range("P200").Value = 571
$ws.Range("Q200").Value = 35
$ws.Range("R200").Value = "Hortaliza"

# --- New row 201 ---
$ws.Range("A201").Value = 10
$ws.Range("B201").Value = "Vega Modelo de Temuco"
$ws.Range("C201").Value = "La Araucanía"
$ws.Range("D201").Value = 45089
$ws.Range("E201").Value = 9
$ws.Range("F201").Value = 100112013
$ws.Range("G201").Value = "Alcachofa"
$ws.Range("H201").Value = "Española"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 400
$ws.Range("K201").Value = 550
$ws.Range("L201").Value = 550
$ws.Range("M201").Value = 550
$ws.Range("N201").Value = "$/unidad"
$ws.Range("O201").Value = "Provincia de Limarí"
$ws.Range("P201").Value = 550
$ws.Range("Q201").Value = 1
$ws.Range("R201").Value = "Hortaliza"
